# Teleportsteine sind vorhanden, aber noch keine Funktion
# Add two new "Teleport stone" rows (TeleA / TeleB) with their RGB values
# to the RGB-Werte table, right after the existing "Weg" row (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "TeleA"
$ws.Range("B13").Value = 250
$ws.Range("C13").Value = 250
$ws.Range("D13").Value = 250

$ws.Range("A14").Value = "TeleB"
$ws.Range("B14").Value = 200
$ws.Range("C14").Value = 200
$ws.Range("D14").Value = 200

# Match the resulting selection state seen in the saved workbook
[void]$ws.Range("D15").Select()
